$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 20 new rows before the old row 48 ("Accounts Payable (Balance Sheet)"),
# shifting all subsequent rows down by 20.
$ws.Rows("48:67").Insert()

# Populate the newly inserted rows 48-67 with the new summary-statistics data
# (Altman Z-score components and Tone-related constructed variables).
$newData = New-Object 'object[,]' 20,7
$newData[0,0] = 'Common Plus Preferred Stock'
$newData[0,1] = 338752830.5819008
$newData[0,2] = -74100
$newData[0,3] = 6000000
$newData[0,4] = 9817134000
$newData[0,5] = 926606262.0503973
$newData[0,6] = 'Constructed for Altman''s Z'
$newData[1,0] = 'EBIT'
$newData[1,1] = 304424183.9851524
$newData[1,2] = -1364004000
$newData[1,3] = 122944000
$newData[1,4] = 4334000000
$newData[1,5] = 556409802.6039112
$newData[1,6] = 'Constructed for Altman''s Z'
$newData[2,0] = 'Ratio A'
$newData[2,1] = 0.02133813934414754
$newData[2,2] = -0.02096513111912044
$newData[2,3] = 0.01912278186664593
$newData[2,4] = 0.07929460716437682
$newData[2,5] = 0.01868069966686935
$newData[2,6] = 'Constructed for Altman''s Z'
$newData[3,0] = 'Ratio B'
$newData[3,1] = 0.2134202507865269
$newData[3,2] = 0.03648170548826261
$newData[3,3] = 0.1774521624681313
$newData[3,4] = 0.7012045786049858
$newData[3,5] = 0.1508570507146134
$newData[3,6] = 'Constructed for Altman''s Z'
$newData[4,0] = 'Ratio C'
$newData[4,1] = 1.868819499357538
$newData[4,2] = 0.2927606442877744
$newData[4,3] = 1.416164556962025
$newData[4,4] = 8.05944049356596
$newData[4,5] = 1.554452326766635
$newData[4,6] = 'Constructed for Altman''s Z'
$newData[5,0] = 'Ratio D'
$newData[5,1] = 0.1328359382574212
$newData[5,2] = -0.1274138657291709
$newData[5,3] = 0.1023823312099367
$newData[5,4] = 0.570361917093548
$newData[5,5] = 0.1511480210306297
$newData[5,6] = 'Constructed for Altman''s Z'
$newData[6,0] = 'Ratio E'
$newData[6,1] = 0.2251411477369084
$newData[6,2] = -0.7448346525137804
$newData[6,3] = 0.2089439925466729
$newData[6,4] = 1.010775408128855
$newData[6,5] = 0.3156237391664442
$newData[6,6] = 'Constructed for Altman''s Z'
$newData[7,0] = 'Working Capital'
$newData[7,1] = 1125108587.771274
$newData[7,2] = -28931855000
$newData[7,3] = 543614000
$newData[7,4] = 39464552600
$newData[7,5] = 3845915891.898036
$newData[7,6] = 'Constructed for Altman''s Z'
$newData[8,0] = 'Active Tone'
$newData[8,1] = 614.4351061898711
$newData[8,2] = 29
$newData[8,3] = 623
$newData[8,4] = 1584
$newData[8,5] = 180.4774582484301
$newData[8,6] = 'Constructed for Tone'
$newData[9,0] = 'Active-Passive Tone Score'
$newData[9,1] = 2.978833305757044
$newData[9,2] = 1.592715231788079
$newData[9,3] = 2.920454545454545
$newData[9,4] = 5.984615384615385
$newData[9,5] = 0.5450171553903713
$newData[9,6] = 'Constructed for Tone'
$newData[10,0] = 'Negative Tone'
$newData[10,1] = 106.2768197495008
$newData[10,2] = 7
$newData[10,3] = 103
$newData[10,4] = 286
$newData[10,5] = 41.5290862580096
$newData[10,6] = 'Constructed for Tone'
$newData[11,0] = 'Overstated Tone'
$newData[11,1] = 380.0214194953712
$newData[11,2] = 27
$newData[11,3] = 382
$newData[11,4] = 932
$newData[11,5] = 116.0017180779686
$newData[11,6] = 'Constructed for Tone'
$newData[12,0] = 'Overstated-Understated Tone Score'
$newData[12,1] = 1.971750632849826
$newData[12,2] = 0.8177777777777778
$newData[12,3] = 1.909090909090909
$newData[12,4] = 5.925925925925926
$newData[12,5] = 0.4716430009099636
$newData[12,6] = 'Constructed for Tone'
$newData[13,0] = 'Passive Tone'
$newData[13,1] = 212.8749319295698
$newData[13,2] = 8
$newData[13,3] = 213
$newData[13,4] = 578
$newData[13,5] = 72.04862095437315
$newData[13,6] = 'Constructed for Tone'
$newData[14,0] = 'Positive Tone'
$newData[14,1] = 333.1058268288256
$newData[14,2] = 21
$newData[14,3] = 335
$newData[14,4] = 994
$newData[14,5] = 104.9695405754105
$newData[14,6] = 'Constructed for Tone'
$newData[15,0] = 'Postivity-Negativity Tone Score'
$newData[15,1] = 3.428857210882757
$newData[15,2] = 1.1
$newData[15,3] = 3.188976377952756
$newData[15,4] = 11.53191489361702
$newData[15,5] = 1.274432146041699
$newData[15,6] = 'Constructed for Tone'
$newData[16,0] = 'Strong Tone'
$newData[16,1] = 715.2392448720276
$newData[16,2] = 40
$newData[16,3] = 708
$newData[16,4] = 2014
$newData[16,5] = 227.5600061618208
$newData[16,6] = 'Constructed for Tone'
$newData[17,0] = 'Strong-Weak Tone Score'
$newData[17,1] = 8.001598122516048
$newData[17,2] = 2.494505494505495
$newData[17,3] = 7.431578947368421
$newData[17,4] = 35.09375
$newData[17,5] = 2.927189613400644
$newData[17,6] = 'Constructed for Tone'
$newData[18,0] = 'Understated Tone'
$newData[18,1] = 202.9088763840988
$newData[18,2] = 17
$newData[18,3] = 203
$newData[18,4] = 565
$newData[18,5] = 73.83442173512447
$newData[18,6] = 'Constructed for Tone'
$newData[19,0] = 'Weak Tone'
$newData[19,1] = 97.34561626429479
$newData[19,2] = 6
$newData[19,3] = 96
$newData[19,4] = 269
$newData[19,5] = 37.24496414477247
$newData[19,6] = 'Constructed for Tone'

$ws.Range("A48:G67").Value = $newData
